$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the two additional date columns (AK, AL)
$ws.Range("AK1").Value = "25-jul"
$ws.Range("AL1").Value = "26-jul"

# New data values for rows 2-18 in columns AK and AL
$ws.Range("AK2").Value = 0
$ws.Range("AL2").Value = 0

$ws.Range("AK3").Value = 11.659527624372734
$ws.Range("AL3").Value = 13.915561918321162

$ws.Range("AK4").Value = 13.720636393633743
$ws.Range("AL4").Value = 14.251428814301386

$ws.Range("AK5").Value = 32.333350072248543
$ws.Range("AL5").Value = 31.199771520604326

$ws.Range("AK6").Value = 0
$ws.Range("AL6").Value = 0

$ws.Range("AK7").Value = 4.7557679868774043
$ws.Range("AL7").Value = 4.4834732748942949

$ws.Range("AK8").Value = 7.7698180044860479
$ws.Range("AL8").Value = 6.5948098977110403

$ws.Range("AK9").Value = 19.041765407060616
$ws.Range("AL9").Value = 20.055897514355209

$ws.Range("AK10").Value = 21.205082518937406
$ws.Range("AL10").Value = 20.226085172726076

$ws.Range("AK11").Value = 16.046546576595595
$ws.Range("AL11").Value = 14.468429411561974

$ws.Range("AK12").Value = 0
$ws.Range("AL12").Value = 0

$ws.Range("AK13").Value = 17.277466146157163
$ws.Range("AL13").Value = 15.842736838675009

$ws.Range("AK14").Value = 0
$ws.Range("AL14").Value = 0

$ws.Range("AK15").Value = 0
$ws.Range("AL15").Value = 0

$ws.Range("AK16").Value = 6.8181059731661264
$ws.Range("AL16").Value = 5.8811808146019091

$ws.Range("AK17").Value = 0
$ws.Range("AL17").Value = 0

$ws.Range("AK18").Value = 0
$ws.Range("AL18").Value = 0

# Update the selected cell to reflect the active cell after the edit
[void]$ws.Range("AM6").Select()
